# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 262 (shifting the existing rows
# 262-364 down to 263-365) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(262).Insert()

$ws.Range("A262").Value = 5
$ws.Range("B262").Value = "Macroferia Regional de Talca"
$ws.Range("C262").Value = "Maule"
$ws.Range("D262").Value = 45009
$ws.Range("E262").Value = 7
$ws.Range("F262").Value = "Fruta"
$ws.Range("G262").Value = 100108
$ws.Range("H262").Value = "Tropicales y subtropicales"
$ws.Range("I262").Value = 100108005
$ws.Range("J262").Value = "Piña"
$ws.Range("K262").Value = "Caramelo"
$ws.Range("L262").Value = "Segunda"
$ws.Range("M262").Value = 230
$ws.Range("N262").Value = 20000
$ws.Range("O262").Value = 20000
$ws.Range("P262").Value = 20000
$ws.Range("Q262").Value = "$/caja 14 unidades"
$ws.Range("R262").Value = "Ecuador"
$ws.Range("S262").Value = 1429
$ws.Range("T262").Value = 14
